$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 409, shifting existing rows 409.. down by one.
$ws.Rows.Item(409).Insert()

# Populate the newly inserted row 409 with the new record (a duplicate of the
# original row 409 entry, but with an updated date and updated prices).
$ws.Cells.Item(409, 1).Value = 4
$ws.Cells.Item(409, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(409, 3).Value = "Los Lagos"
$ws.Cells.Item(409, 4).Value = 44995
$ws.Cells.Item(409, 5).Value = 10
$ws.Cells.Item(409, 6).Value = 100112023
$ws.Cells.Item(409, 7).Value = "Brócoli"
$ws.Cells.Item(409, 8).Value = "Sin especificar"
$ws.Cells.Item(409, 9).Value = "Primera"
$ws.Cells.Item(409, 10).Value = 500
$ws.Cells.Item(409, 11).Value = 1600
$ws.Cells.Item(409, 12).Value = 1600
$ws.Cells.Item(409, 13).Value = 1600
$ws.Cells.Item(409, 14).Value = "$/unidad"
$ws.Cells.Item(409, 15).Value = "Región Metropolitana"
$ws.Cells.Item(409, 16).Value = 1600
$ws.Cells.Item(409, 17).Value = 1
$ws.Cells.Item(409, 18).Value = "Hortaliza"
